$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 0.530292749404907
$ws.Range("I2").Value = 0.49904227256774902
$ws.Range("H3").Value = 0.66603994369506803
$ws.Range("I3").Value = 0.509785175323486
$ws.Range("H4").Value = 0.67190027236938399
$ws.Range("I4").Value = 0.56935691833496005
$ws.Range("H5").Value = 0.688501596450805
$ws.Range("H6").Value = 14.4302310943603
$ws.Range("H7").Value = 10.007211923599201
$ws.Range("I7").Value = 14.7417678833007
$ws.Range("H8").Value = 13.1977627277374
$ws.Range("H9").Value = 9.8919746875762904
$ws.Range("H10").Value = 5.9474902153015101
$ws.Range("I10").Value = 2.2901258468627899
$ws.Range("H11").Value = 4.3761410713195801
$ws.Range("H12").Value = 8.8802170753479004
$ws.Range("H13").Value = 7.9954185485839799
$ws.Range("I13").Value = 2.08503866195678
$ws.Range("H14").Value = 19.328854560852001
$ws.Range("I14").Value = 1.8281934261321999
$ws.Range("H15").Value = 11.078542470932
$ws.Range("H16").Value = 17.684261083602902
$ws.Range("I16").Value = 3.13488340377807
$ws.Range("H17").Value = 4.5148169994354204
$ws.Range("I17").Value = 6.7404875755309996

$ws.Range("H18").Select()
